# "Added New Mac-Address and Document Types"
# Append 5 new test-data rows (157-161) to the master-reg_center_machine_device
# sheet, following the same pattern as the existing rows immediately above
# them (lang_code "eng", is_active TRUE, cr_by "superadmin", cr_dtimes "now()").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id, machine_id, device_id for each new row
$newRows = @(
    @(10002, 10032, 3000176),
    @(10002, 10032, 3000177),
    @(10002, 10032, 3000178),
    @(10002, 10032, 3000179),
    @(10002, 10032, 3000180)
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]          # regcntr_id
    $ws.Cells.Item($r, 2).Value = $vals[1]          # machine_id
    $ws.Cells.Item($r, 3).Value = $vals[2]          # device_id
    $ws.Cells.Item($r, 4).Value = "eng"             # lang_code
    $ws.Cells.Item($r, 5).Value = $true             # is_active
    $ws.Cells.Item($r, 6).Value = "superadmin"      # cr_by
    $ws.Cells.Item($r, 7).Value = "now()"           # cr_dtimes
}

# Match the post-edit view state: scrolled down with E157 as the active cell.
$ws.Range("E157").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 150
$excel.ActiveWindow.ScrollColumn = 1
